$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Quantity Per Unit" column (column L) entirely; everything to
# the right shifts one column to the left, and the unused shared string is
# dropped automatically.
$ws.Range("L:L").EntireColumn.Delete()

# Restore the selection state recorded in the saved workbook: the user had
# just selected the (now) entire column L ("Item Group").
$ws.Range("L1:L1048576").Select()
